# add free proxy maker
# Append two new user records (rows 20 and 21) to the "Full Data" sheet,
# mirroring how the rows above were filled in (Last/First/Birth/Place/
# Mother/Father/Email/Passport columns A-H).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 20: Oya Evran -------------------------------------------------
$ws.Range("A20").Value = "Evran "
$ws.Range("B20").Value = "Oya "
$ws.Range("E20").Value = "Latife "
$ws.Range("F20").Value = "Ahmet "
$ws.Range("D20").Value = "Vize"
$ws.Range("H20").Value = "D12790130"
$ws.Range("G20").Value = "birsenaltan1@hotmail.com"

# --- Row 21: MuhammedSalih Karademir -----------------------------------
$ws.Range("A21").Value = "Karademir"
$ws.Range("B21").Value = "MuhammedSalih"
$ws.Range("D21").Value = "Sakarya"
$ws.Range("E21").Value = "Selver"
$ws.Range("H21").Value = "U31185115"
$ws.Range("G21").Value = "salihkarademir17@gmail.com"
$ws.Range("F21").Value = "Yuksel"

# Birth dates are entered last, as text (one of them - 1997-28-12 - isn't
# even a valid calendar date), so the column is formatted as Text first.
$ws.Range("C20:C21").NumberFormat = "@"
$ws.Range("C20").Value = "1969-04-30"
$ws.Range("C21").Value = "1997-28-12"

# Re-fit the first three columns now that they hold longer values.
$ws.Range("A1:C21").EntireColumn.AutoFit()

# Leave the selection where the user finished typing.
[void]$ws.Range("G20").Select()

# Match the printed page setup recorded for this sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
